$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.3260275
$ws.Range("H2").Value = 0.6520550000000001
$ws.Range("I2").Value = 0.4722113996121241
$ws.Range("J2").Value = 0.4126724043544658
$ws.Range("O2").Value = 0.367191313355536
$ws.Range("P2").Value = 0.368867471188671
$ws.Range("Q2").Value = 0.7466986097333334
$ws.Range("R2").Value = 4.4801916584
$ws.Range("S2").Value = 0.1733919240050317
$ws.Range("T2").Value = 0.1522214262235805
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.3260275
$ws.Range("H3").Value = 0.6520550000000001
$ws.Range("I3").Value = 0.4722113996121241
$ws.Range("J3").Value = 0.4126724043544658
$ws.Range("M3").Value = 1.867716
$ws.Range("N3").Value = 5.603148000000001
$ws.Range("O3").Value = 0.2994415959884972
$ws.Range("P3").Value = 0.3008084893719378
$ws.Range("Q3").Value = 0.6089267781900002
$ws.Range("R3").Value = 3.653560669140001
$ws.Range("S3").Value = 0.1413997351438165
$ws.Range("T3").Value = 0.1241353625593524
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.3260275
$ws.Range("H4").Value = 0.6520550000000001
$ws.Range("I4").Value = 0.4722113996121241
$ws.Range("J4").Value = 0.4126724043544658
$ws.Range("M4").Value = 0.08502850000000001
$ws.Range("N4").Value = 0.170057
$ws.Range("O4").Value = 0.01363219555034488
$ws.Range("P4").Value = 0.00912961593681331
$ws.Range("Q4").Value = 0.02772162928375
$ws.Range("R4").Value = 0.110886517135
$ws.Range("S4").Value = 0.006437278140614524
$ws.Range("T4").Value = 0.003767540559477598
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.3260275
$ws.Range("H5").Value = 0.6520550000000001
$ws.Range("I5").Value = 0.4722113996121241
$ws.Range("J5").Value = 0.4126724043544658
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.132021
$ws.Range("N5").Value = 3.396063
$ws.Range("O5").Value = 0.1814912839706329
$ws.Range("P5").Value = 0.1823197568298983
$ws.Range("Q5").Value = 0.3690699765775
$ws.Range("R5").Value = 2.214419859465
$ws.Range("S5").Value = 0.085702253221174
$ws.Range("T5").Value = 0.07523833241231569
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3260275
$ws.Range("H6").Value = 0.6520550000000001
$ws.Range("I6").Value = 0.4722113996121241
$ws.Range("J6").Value = 0.4126724043544658
$ws.Range("M6").Value = 0.8622709999999999
$ws.Range("N6").Value = 2.586813
$ws.Range("O6").Value = 0.1382436111349891
$ws.Range("P6").Value = 0.1388746666726794
$ws.Range("Q6").Value = 0.2811240584525
$ws.Range("R6").Value = 1.686744350715
$ws.Range("S6").Value = 0.06528020910148745
$ws.Range("T6").Value = 0.05730974259973964
$ws.Range("I7").Value = 0.2885533155568945
$ws.Range("J7").Value = 0.3782564036355144
$ws.Range("O7").Value = 0.367191313355536
$ws.Range("P7").Value = 0.368867471188671
$ws.Range("S7").Value = 0.1059542709124305
$ws.Range("T7").Value = 0.1395264830699534
$ws.Range("I8").Value = 0.2885533155568945
$ws.Range("J8").Value = 0.3782564036355144
$ws.Range("M8").Value = 1.867716
$ws.Range("N8").Value = 5.603148000000001
$ws.Range("O8").Value = 0.2994415959884972
$ws.Range("P8").Value = 0.3008084893719378
$ws.Range("Q8").Value = 0.3720957201
$ws.Range("R8").Value = 3.3488614809
$ws.Range("S8").Value = 0.08640486533812893
$ws.Range("T8").Value = 0.1137827373728611
$ws.Range("I9").Value = 0.2885533155568945
$ws.Range("J9").Value = 0.3782564036355144
$ws.Range("M9").Value = 0.08502850000000001
$ws.Range("N9").Value = 0.170057
$ws.Range("O9").Value = 0.01363219555034488
$ws.Range("P9").Value = 0.00912961593681331
$ws.Range("Q9").Value = 0.0169398029125
$ws.Range("R9").Value = 0.101638817475
$ws.Range("S9").Value = 0.003933615224371957
$ws.Range("T9").Value = 0.00345333569083248
$ws.Range("I10").Value = 0.2885533155568945
$ws.Range("J10").Value = 0.3782564036355144
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.132021
$ws.Range("N10").Value = 3.396063
$ws.Range("O10").Value = 0.1814912839706329
$ws.Range("P10").Value = 0.1823197568298983
$ws.Range("Q10").Value = 0.225526883725
$ws.Range("R10").Value = 2.029741953525
$ws.Range("S10").Value = 0.05236991173440396
$ws.Range("T10").Value = 0.06896361553017885
$ws.Range("I11").Value = 0.2885533155568945
$ws.Range("J11").Value = 0.3782564036355144
$ws.Range("M11").Value = 0.8622709999999999
$ws.Range("N11").Value = 2.586813
$ws.Range("O11").Value = 0.1382436111349891
$ws.Range("P11").Value = 0.1388746666726794
$ws.Range("Q11").Value = 0.171785939975
$ws.Range("R11").Value = 1.546073459775
$ws.Range("S11").Value = 0.03989065234755913
$ws.Range("T11").Value = 0.05253023197168855
$ws.Range("I12").Value = 0.2392352848309814
$ws.Range("J12").Value = 0.2090711920100198
$ws.Range("O12").Value = 0.367191313355536
$ws.Range("P12").Value = 0.368867471188671
$ws.Range("S12").Value = 0.08784511843807381
$ws.Range("T12").Value = 0.07711956189513706
$ws.Range("I13").Value = 0.2392352848309814
$ws.Range("J13").Value = 0.2090711920100198
$ws.Range("M13").Value = 1.867716
$ws.Range("N13").Value = 5.603148000000001
$ws.Range("O13").Value = 0.2994415959884972
$ws.Range("P13").Value = 0.3008084893719378
$ws.Range("Q13").Value = 0.308499056442
$ws.Range("R13").Value = 1.850994338652
$ws.Range("S13").Value = 0.07163699550655178
$ws.Range("T13").Value = 0.0628903894397244
$ws.Range("I14").Value = 0.2392352848309814
$ws.Range("J14").Value = 0.2090711920100198
$ws.Range("M14").Value = 0.08502850000000001
$ws.Range("N14").Value = 0.170057
$ws.Range("O14").Value = 0.01363219555034488
$ws.Range("P14").Value = 0.00912961593681331
$ws.Range("Q14").Value = 0.01404453997325
$ws.Range("R14").Value = 0.05617815989300001
$ws.Range("S14").Value = 0.003261302185358394
$ws.Range("T14").Value = 0.001908739686503232
$ws.Range("I15").Value = 0.2392352848309814
$ws.Range("J15").Value = 0.2090711920100198
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1.132021
$ws.Range("N15").Value = 3.396063
$ws.Range("O15").Value = 0.1814912839706329
$ws.Range("P15").Value = 0.1823197568298983
$ws.Range("Q15").Value = 0.1869810026645
$ws.Range("R15").Value = 1.121886015987
$ws.Range("S15").Value = 0.04341911901505488
$ws.Range("T15").Value = 0.03811780888740378
$ws.Range("I16").Value = 0.2392352848309814
$ws.Range("J16").Value = 0.2090711920100198
$ws.Range("M16").Value = 0.8622709999999999
$ws.Range("N16").Value = 2.586813
$ws.Range("O16").Value = 0.1382436111349891
$ws.Range("P16").Value = 0.1388746666726794
$ws.Range("Q16").Value = 0.1424251812895
$ws.Range("R16").Value = 0.8545510877369999
$ws.Range("S16").Value = 0.03307274968594257
$ws.Range("T16").Value = 0.02903469210125125

Write-Host "Applied $([int]166) cell updates to sheet '$($ws.Name)'."
